$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.589.56"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +4.69%  "

$ws.Range("D3").Value = "'2.722.52"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.37%  "

$ws.Range("E4").Value = "  +0.18%  "

$ws.Range("D5").Value = "'576.62"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.06%  "

$ws.Range("D6").Value = "'154.10"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.25%  "

$ws.Range("D7").Value = "'0.996"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.18%  "

$ws.Range("E8").Value = "  +1.48%  "

$ws.Range("D9").Value = "'2.747.41"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.93%  "

$ws.Range("D10").Value = "'6.68"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.83%  "

$ws.Range("E11").Value = "  +5.37%  "

$ws.Range("E12").Value = "  +4.92%  "

$ws.Range("E13").Value = "  +3.77%  "

$ws.Range("D14").Value = "'3.208.11"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.50%  "

$ws.Range("D15").Value = "'26.39"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.57%  "

$ws.Range("D16").Value = "'63.640.61"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.81%  "

$ws.Range("E17").Value = "  +6.47%  "

$ws.Range("D18").Value = "'2.742.75"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.08%  "

$ws.Range("D19").Value = "'11.91"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.14%  "

$ws.Range("D20").Value = "'4.86"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.70%  "

$ws.Range("D21").Value = "'359.90"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.91%  "

$ws.Range("D22").Value = "'6.92"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.32%  "

$ws.Range("D23").Value = "'0.996"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.30%  "

$ws.Range("D24").Value = "'0.534"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.11%  "

$ws.Range("D25").Value = "'65.78"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.82%  "

$ws.Range("D26").Value = "'0.169"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.48%  "

$ws.Range("E27").Value = "  +5.35%  "

$ws.Range("D28").Value = "'0.999"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.05%  "

$ws.Range("D29").Value = "'0.0₃0903"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +12.40%  "

$ws.Range("E30").Value = "  -1.19%  "

$ws.Range("D31").Value = "'7.13"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +7.09%  "

$ws.Range("D32").Value = "'171.72"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.58%  "

$ws.Range("E33").Value = "  +13.80%  "

$ws.Range("D34").Value = "'0.997"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.07%  "

$ws.Range("D35").Value = "'20.46"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.18%  "

$ws.Range("E36").Value = "  +6.92%  "

$ws.Range("E37").Value = "  +9.19%  "

$ws.Range("E38").Value = "  +10.07%  "

$ws.Range("D39").Value = "'0.998"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +14.38%  "

$ws.Range("D40").Value = "'343.39"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.01%  "

$ws.Range("D41").Value = "'4.21"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.95%  "

$ws.Range("D42").Value = "'39.24"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.22%  "

$ws.Range("D43").Value = "'5.54"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +6.41%  "

$ws.Range("D44").Value = "'21.70"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +7.78%  "

$ws.Range("D45").Value = "'21.82"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.57%  "

$ws.Range("D46").Value = "'0.0590"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.66%  "

$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "'138.95"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.05%  "

$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").Value = "'0.644"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +5.18%  "

$ws.Range("D49").Value = "'0.0254"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.02%  "

$ws.Range("E50").Value = "  +1.05%  "

$ws.Range("D51").Value = "'0.996"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.30%  "
